$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Comments" header in column J and a comment for the last row.
$ws.Range("J1").Value = "Comments"
$ws.Range("J11").Value = "This has never been cultured so has an unusual name"

# Style the new comment cell: centered, wrapped text.
$ws.Range("J11").HorizontalAlignment = -4108
$ws.Range("J11").WrapText = $true

# Widen column J to fit the new comments text and let row 11 grow to fit the wrapped text.
$ws.Columns.Item(10).ColumnWidth = 22.5
$ws.Rows.Item(11).RowHeight = 32

# Move the active selection, matching the author's final cursor position.
$ws.Range("H22").Select()
